$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new values are plain decimal numbers (e.g. "1.00", "0.746").
# Assigning them straight to .Value would make Excel auto-convert the
# text into a real number (stripping the original text formatting, like
# "1.00" -> 1). Pin those cells to Text format before the write, then
# restore the default style afterwards so no extra formatting lingers.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D15", "D18", "D21", "D22", "D24", "D26", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D40", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.017.25"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.860.86"
$ws.Range("D5").Value = "472.84"
$ws.Range("E5").Value = "  +10.22%  "
$ws.Range("D6").Value = "144.84"
$ws.Range("E6").Value = "  +10.14%  "
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "0.746"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").Value = "0.0000312"
$ws.Range("E11").Value = "  -7.48%  "
$ws.Range("D12").Value = "43.48"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").Value = "10.43"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "4.483.54"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").Value = "14.83"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").Value = "3.890.51"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "20.11"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("D20").Value = "67.303.94"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "431.82"
$ws.Range("E21").Value = "  +3.79%  "
$ws.Range("D22").Value = "15.01"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  +6.13%  "
$ws.Range("D24").Value = "88.44"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("E25").Value = "  +8.79%  "
$ws.Range("D26").Value = "37.98"
$ws.Range("E27").Value = "  +6.70%  "
$ws.Range("D28").Value = "9.95"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "5.54"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").Value = "727.92"
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").Value = "13.93"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("E32").Value = "  +6.64%  "
$ws.Range("D33").Value = "2.76"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "43.36"
$ws.Range("E34").Value = "  +11.23%  "
$ws.Range("D35").Value = "0.161"
$ws.Range("E35").Value = "  +7.69%  "
$ws.Range("D36").Value = "58.43"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  -6.12%  "
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").Value = "0.348"
$ws.Range("E40").Value = "  +7.40%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").Value = "  +2.97%  "
$ws.Range("D43").Value = "0.0₃0677"
$ws.Range("E43").Value = "  -7.10%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  +5.26%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "3.47"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("E47").Value = "  +5.38%  "
$ws.Range("E48").Value = "  +5.10%  "
$ws.Range("D49").Value = "3.19"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "2.91"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").Value = "143.37"
$ws.Range("E51").Value = "  +1.35%  "

# Restore default (unstyled) formatting on the cells we pinned to Text,
# now that the text value is safely committed.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
